# Forms the consolidated report: recompute the "Absent" column (H) from the
# "Real" attendance column (E) for each date row in the attendance sheet.
# A student is considered Absent (H = 1) on a given day when their Real
# attendance count (E) is 0; otherwise they are not absent (H = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $realValue = $ws.Cells.Item($r, 5).Value2
    if ($realValue -eq 0) {
        $absent = 1
    } else {
        $absent = 0
    }
    $ws.Cells.Item($r, 8).Value = $absent
}
